$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6038436889648438
$ws.Range("B1").Value = 1.507175803184509
$ws.Range("C1").Value = 5.873884677886963
$ws.Range("D1").Value = 2.10520339012146
$ws.Range("E1").Value = 1.465823769569397
